# Applies the scheduled-runner profit-data refresh described in the commit.
# For each affected sheet, numeric columns H-N (currentAveragePrice..LeveProfitHQ)
# are updated to the newly fetched market-board figures. A handful of cells that
# no longer carry a value (formerly 0/derived) are cleared outright rather than
# written as 0, matching how the upstream exporter omits empty numeric cells.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 6163.6665
$ws.Range("J40").Value = 5188.846
$ws.Range("L40").Value = 5188.846
$ws.Range("N40").Value = -5538.846
$ws.Range("H64").Value = 5796.885
$ws.Range("I64").Value = 5316.2104
$ws.Range("J64").Value = 7101.5713
$ws.Range("K64").Value = 5316.2104
$ws.Range("L64").Value = 7101.5713
$ws.Range("M64").Value = -5068.2104
$ws.Range("N64").Value = -7597.5713
$ws.Range("H67").Value = 5796.885
$ws.Range("I67").Value = 5316.2104
$ws.Range("J67").Value = 7101.5713
$ws.Range("K67").Value = 5316.2104
$ws.Range("L67").Value = 7101.5713
$ws.Range("M67").Value = -4458.2104
$ws.Range("N67").Value = -8817.5713
$ws.Range("H103").Value = 1026
$ws.Range("I103").Value = 1373.25
$ws.Range("J103").Value = 919.1539
$ws.Range("K103").Value = 4119.75
$ws.Range("L103").Value = 2757.4617
$ws.Range("M103").Value = -3533.75
$ws.Range("N103").Value = -3929.4617
$ws.Range("H125").Value = 3008.7273
$ws.Range("I125").Value = 1959
$ws.Range("K125").Value = 17631
$ws.Range("M125").Value = -15171
$ws.Range("H137").Value = 11743.083
$ws.Range("I137").Value = 1789.5714
$ws.Range("J137").Value = 25678
$ws.Range("K137").Value = 5368.7142
$ws.Range("L137").Value = 77034
$ws.Range("M137").Value = -2818.7142
$ws.Range("N137").Value = -82134

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 13892104
$ws.Range("I32").Value = 15876104
$ws.Range("K32").Value = 15876104
$ws.Range("M32").Value = -15875817
$ws.Range("H61").Value = 40001572
$ws.Range("I61").Value = 52633124
$ws.Range("K61").Value = 52633124
$ws.Range("M61").Value = -52632912
$ws.Range("H74").Value = 19233072
$ws.Range("I74").Value = 41668280
$ws.Range("K74").Value = 41668280
$ws.Range("M74").Value = -41667406
$ws.Range("H77").Value = 19233072
$ws.Range("I77").Value = 41668280
$ws.Range("K77").Value = 208341400
$ws.Range("M77").Value = -208337032
$ws.Range("H136").Value = 40001572
$ws.Range("I136").Value = 52633124
$ws.Range("K136").Value = 157899372
$ws.Range("M136").Value = -157896822

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H74").Value = 40387.75
$ws.Range("I74").Value = 49994
$ws.Range("J74").Value = 37185.668
$ws.Range("K74").Value = 49994
$ws.Range("L74").Value = 37185.668
$ws.Range("M74").Value = -49058
$ws.Range("N74").Value = -39057.668
$ws.Range("H77").Value = 40387.75
$ws.Range("I77").Value = 49994
$ws.Range("J77").Value = 37185.668
$ws.Range("K77").Value = 149982
$ws.Range("L77").Value = 111557.004
$ws.Range("M77").Value = -145302
$ws.Range("N77").Value = -120917.004
$ws.Range("H82").Value = 42138.832
$ws.Range("I82").Value = 17700
$ws.Range("J82").Value = 47026.6
$ws.Range("K82").Value = 17700
$ws.Range("L82").Value = 47026.6
$ws.Range("M82").Value = -17317
$ws.Range("N82").Value = -47792.6
$ws.Range("H85").Value = 42138.832
$ws.Range("I85").Value = 17700
$ws.Range("J85").Value = 47026.6
$ws.Range("K85").Value = 17700
$ws.Range("L85").Value = 47026.6
$ws.Range("M85").Value = -16374
$ws.Range("N85").Value = -49678.6

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H17").Value = 15013
$ws.Range("J17").Value = 15013
$ws.Range("L17").Value = 15013
$ws.Range("N17").Value = -15361
$ws.Range("H31").Value = 59529284
$ws.Range("I31").Value = 3739.1538
$ws.Range("J31").Value = 156258300
$ws.Range("K31").Value = 3739.1538
$ws.Range("L31").Value = 156258300
$ws.Range("M31").Value = -3444.1538
$ws.Range("N31").Value = -156258890
$ws.Range("H34").Value = 59529284
$ws.Range("I34").Value = 3739.1538
$ws.Range("J34").Value = 156258300
$ws.Range("K34").Value = 3739.1538
$ws.Range("L34").Value = 156258300
$ws.Range("M34").Value = -3537.1538
$ws.Range("N34").Value = -156258704
$ws.Range("H107").Value = 536.62964
$ws.Range("I107").Value = 473.5909
$ws.Range("K107").Value = 473.5909
$ws.Range("M107").Value = 1446.4091
$ws.Range("I132").Value = 2712.2888
$ws.Range("J132").Value = 3162.2222
$ws.Range("K132").Value = 8136.866399999999
$ws.Range("L132").Value = 9486.6666
$ws.Range("M132").Value = -5606.866399999999
$ws.Range("N132").Value = -14546.6666

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H98").Value = 1846
$ws.Range("I98").Value = 2692.5
$ws.Range("J98").Value = 999.5
$ws.Range("K98").Value = 8077.5
$ws.Range("L98").Value = 2998.5
$ws.Range("M98").Value = -6579.5
$ws.Range("N98").Value = -5994.5
$ws.Range("H129").Value = 27779268
$ws.Range("J129").Value = 20835120
$ws.Range("L129").Value = 62505360
$ws.Range("N129").Value = -62515360
$ws.Range("H140").Value = 1988.2354
$ws.Range("I140").Value = 1731
$ws.Range("J140").Value = 3188.6667
$ws.Range("K140").Value = 5193
$ws.Range("L140").Value = 9566.000100000001
$ws.Range("M140").Value = -13
$ws.Range("N140").Value = -19926.0001

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 3300
$ws.Range("I22").Value = 2916.6667
$ws.Range("J22").Value = 3875
$ws.Range("K22").Value = 2916.6667
$ws.Range("L22").Value = 3875
$ws.Range("M22").Value = -2621.6667
$ws.Range("N22").Value = -4465
$ws.Range("H27").Value = 3300
$ws.Range("I27").Value = 2916.6667
$ws.Range("J27").Value = 3875
$ws.Range("K27").Value = 2916.6667
$ws.Range("L27").Value = 3875
$ws.Range("M27").Value = -2809.6667
$ws.Range("N27").Value = -4089
$ws.Range("H55").Value = 536
$ws.Range("I55").Value = 493.4
$ws.Range("J55").Value = 589.25
$ws.Range("K55").Value = 493.4
$ws.Range("L55").Value = 589.25
$ws.Range("M55").Value = -320.4
$ws.Range("N55").Value = -935.25
$ws.Range("H136").Value = 1179442.4
$ws.Range("I136").Value = 1336034.8
$ws.Range("K136").Value = 4008104.4
$ws.Range("M136").Value = -4005554.4

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H22").Value = 15000
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 20000
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = -9707
$ws.Range("N22").Value = -20586
$ws.Range("H24").Value = 9950
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 28800
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H45").Value = 23666
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 23666
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 23666
$ws.Range("N45").Value = -24648
$ws.Range("M45").ClearContents()
$ws.Range("H96").Value = 4202.4443
$ws.Range("J96").Value = 2599.5
$ws.Range("L96").Value = 2599.5
$ws.Range("N96").Value = -5345.5
$ws.Range("H107").Value = 2117.1177
$ws.Range("J107").Value = 2474.25
$ws.Range("L107").Value = 7422.75
$ws.Range("N107").Value = -11262.75
$ws.Range("H111").Value = 72000
$ws.Range("J111").Value = 72000
$ws.Range("L111").Value = 72000
$ws.Range("N111").Value = -80180
$ws.Range("H132").Value = 1269.88
$ws.Range("I132").Value = 1157.75
$ws.Range("J132").Value = 1718.4
$ws.Range("K132").Value = 3473.25
$ws.Range("L132").Value = 5155.200000000001
$ws.Range("M132").Value = -943.25
$ws.Range("N132").Value = -10215.2
$ws.Range("H136").Value = 1425.8368
$ws.Range("I136").Value = 1450.9778
$ws.Range("J136").Value = 1143
$ws.Range("K136").Value = 4352.9334
$ws.Range("L136").Value = 3429
$ws.Range("M136").Value = -1802.9334
